$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Marking" row (row 11): Right = 4 (was 5), Wrong = -2 (was -1)
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# "Total" row (row 12): Right = 80 (was 100), Wrong = -6 (was -3)
$ws.Range("B12").Value = 80
$ws.Range("C12").Value = -6

# Max column summary text for the Total row
$ws.Range("E12").Value = "74 / 112"
